$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 39:40. This pushes the existing rows 39-56
# (the weekly Haba price records) down to rows 41-58, matching the diff
# where the old row 39 data now lives at row 41, old row 40 at row 42, etc.
$ws.Rows("39:40").Insert()

# Populate the two brand-new weekly records at rows 39 and 40.
$ws.Range("A39").Value = 3
$ws.Range("B39").Value = "Femacal de La Calera"
$ws.Range("C39").Value = "Coquimbo"
$ws.Range("D39").Value = 44455
$ws.Range("E39").Value = 5
$ws.Range("F39").Value = 100112026
$ws.Range("G39").Value = "Haba"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 70
$ws.Range("K39").Value = 11000
$ws.Range("L39").Value = 12000
$ws.Range("M39").Value = 11500
$ws.Range("N39").Value = "$/saco 25 kilos"
$ws.Range("O39").Value = "Provincia de Limarí"
$ws.Range("P39").Value = 460
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"

$ws.Range("A40").Value = 3
$ws.Range("B40").Value = "Femacal de La Calera"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44455
$ws.Range("E40").Value = 5
$ws.Range("F40").Value = 100112026
$ws.Range("G40").Value = "Haba"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Segunda"
$ws.Range("J40").Value = 38
$ws.Range("K40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = 10000
$ws.Range("N40").Value = "$/saco 25 kilos"
$ws.Range("O40").Value = "Provincia de Limarí"
$ws.Range("P40").Value = 400
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = "Hortaliza"
